$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.517.93'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.729.47'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.09'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4828'
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2666'
$ws.Range("E8").Value = '  -0.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06236'
$ws.Range("E9").Value = '  -0.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.730.33'
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07078'
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.570'
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6087'
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.26'
$ws.Range("E15").Value = '  -0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.521.23'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("B18").Value = 'BinanceUSD'
$ws.Range("C18").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007281'
$ws.Range("E19").Value = '  +5.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.52'
$ws.Range("E20").Value = '  -1.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.952.28'
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.499'
$ws.Range("E22").Value = '  -2.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.767'
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.242'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '137.09'
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.41'
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  -2.43%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.406'
$ws.Range("E28").Value = '  -1.80%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '108.25'
$ws.Range("E29").Value = '  +1.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.964'
$ws.Range("E30").Value = '  -1.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08005'
$ws.Range("E31").Value = '  +1.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.690'
$ws.Range("E32").Value = '  -1.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04575'
$ws.Range("E33").Value = '  -0.64%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.002'
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6332'
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.8905'
$ws.Range("E38").Value = '  -5.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.011'
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.401'
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.002'
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01503'
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.83'
$ws.Range("E43").Value = '  -10.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.494'
$ws.Range("E44").Value = '  -4.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3888'
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.001'
$ws.Range("E46").Value = '  +4.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1183'
$ws.Range("E47").Value = '  -2.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05392'
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.918'
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.63'
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.252'
$ws.Range("E51").Value = '  -1.41%  '
